$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (Coin name / Link / Volume columns) - Excel will not
# misinterpret these as numbers, so a direct .Value assignment is safe.
$textUpdates = @(
    @{ Cell = 'E2'; Value = '  +2.05%  ' },
    @{ Cell = 'E3'; Value = '  +2.14%  ' },
    @{ Cell = 'E4'; Value = '  -0.14%  ' },
    @{ Cell = 'E5'; Value = '  +0.37%  ' },
    @{ Cell = 'E6'; Value = '  -0.04%  ' },
    @{ Cell = 'E7'; Value = '  -0.20%  ' },
    @{ Cell = 'E8'; Value = '  +1.19%  ' },
    @{ Cell = 'E9'; Value = '  +1.10%  ' },
    @{ Cell = 'E10'; Value = '  +1.31%  ' },
    @{ Cell = 'E11'; Value = '  +3.41%  ' },
    @{ Cell = 'E12'; Value = '  +1.83%  ' },
    @{ Cell = 'E13'; Value = '  -0.99%  ' },
    @{ Cell = 'E14'; Value = '  +2.34%  ' },
    @{ Cell = 'E15'; Value = '  +1.93%  ' },
    @{ Cell = 'E16'; Value = '  +0.31%  ' },
    @{ Cell = 'E17'; Value = '  +0.71%  ' },
    @{ Cell = 'E18'; Value = '  -0.09%  ' },
    @{ Cell = 'E19'; Value = '  +0.76%  ' },
    @{ Cell = 'E21'; Value = '  -0.08%  ' },
    @{ Cell = 'E22'; Value = '  +2.05%  ' },
    @{ Cell = 'E23'; Value = '  +1.77%  ' },
    @{ Cell = 'E24'; Value = '  +0.81%  ' },
    @{ Cell = 'E25'; Value = '  +0.19%  ' },
    @{ Cell = 'E26'; Value = '  -2.71%  ' },
    @{ Cell = 'E28'; Value = '  +1.97%  ' },
    @{ Cell = 'E29'; Value = '  +1.61%  ' },
    @{ Cell = 'E30'; Value = '  +1.13%  ' },
    @{ Cell = 'E31'; Value = '  +0.37%  ' },
    @{ Cell = 'E32'; Value = '  +0.25%  ' },
    @{ Cell = 'E33'; Value = '  +0.79%  ' },
    @{ Cell = 'E34'; Value = '  +2.04%  ' },
    @{ Cell = 'E35'; Value = '  +0.69%  ' },
    @{ Cell = 'E36'; Value = '  -1.82%  ' },
    @{ Cell = 'E37'; Value = '  +0.83%  ' },
    @{ Cell = 'E38'; Value = '  +1.21%  ' },
    @{ Cell = 'E39'; Value = '  +3.67%  ' },
    @{ Cell = 'E40'; Value = '  +0.33%  ' },
    @{ Cell = 'E41'; Value = '  +2.71%  ' },
    @{ Cell = 'E42'; Value = '  +1.49%  ' },
    @{ Cell = 'E43'; Value = '  +2.99%  ' },
    @{ Cell = 'E44'; Value = '  +3.01%  ' },
    @{ Cell = 'E45'; Value = '  +8.39%  ' },
    @{ Cell = 'E46'; Value = '  +2.44%  ' },
    @{ Cell = 'E47'; Value = '  -1.40%  ' },
    @{ Cell = 'E48'; Value = '  +3.67%  ' },
    @{ Cell = 'E49'; Value = '  +6.70%  ' },
    @{ Cell = 'B50'; Value = 'Quant' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' },
    @{ Cell = 'E50'; Value = '  -0.41%  ' },
    @{ Cell = 'B51'; Value = 'WEMIXToken' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'E51'; Value = '  -4.41%  ' }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Price column updates. These look numeric (e.g. "1.003", "327.80") but the
# source data stores them as literal text (note trailing zeros / multi-dot
# thousand separators like "29.201.37" that would be mangled as a real
# number). Force text interpretation via NumberFormat "@" before writing,
# then restore the "Normal" cell style so no stray formatting is left behind.
$priceUpdates = @(
    @{ Cell = 'D2'; Value = '29.201.37' },
    @{ Cell = 'D3'; Value = '1.910.09' },
    @{ Cell = 'D4'; Value = '1.003' },
    @{ Cell = 'D5'; Value = '327.80' },
    @{ Cell = 'D7'; Value = '0.4645' },
    @{ Cell = 'D8'; Value = '0.3930' },
    @{ Cell = 'D9'; Value = '46.97' },
    @{ Cell = 'D10'; Value = '0.07972' },
    @{ Cell = 'D12'; Value = '22.36' },
    @{ Cell = 'D13'; Value = '1.868.71' },
    @{ Cell = 'D14'; Value = '7.150' },
    @{ Cell = 'D15'; Value = '5.804' },
    @{ Cell = 'D16'; Value = '0.06999' },
    @{ Cell = 'D17'; Value = '88.68' },
    @{ Cell = 'D19'; Value = '0.00001012' },
    @{ Cell = 'D20'; Value = '17.28' },
    @{ Cell = 'D21'; Value = '1.004' },
    @{ Cell = 'D22'; Value = '29.204.77' },
    @{ Cell = 'D24'; Value = '11.08' },
    @{ Cell = 'D25'; Value = '2.112.38' },
    @{ Cell = 'D26'; Value = '2.057' },
    @{ Cell = 'D27'; Value = '155.84' },
    @{ Cell = 'D29'; Value = '5.860' },
    @{ Cell = 'D30'; Value = '2.008' },
    @{ Cell = 'D32'; Value = '0.09392' },
    @{ Cell = 'D33'; Value = '0.9263' },
    @{ Cell = 'D34'; Value = '5.376' },
    @{ Cell = 'D35'; Value = '1.347' },
    @{ Cell = 'D38'; Value = '1.159' },
    @{ Cell = 'D39'; Value = '8.026' },
    @{ Cell = 'D41'; Value = '0.5769' },
    @{ Cell = 'D43'; Value = '10.02' },
    @{ Cell = 'D44'; Value = '12.06' },
    @{ Cell = 'D45'; Value = '2.253' },
    @{ Cell = 'D46'; Value = '0.5436' },
    @{ Cell = 'D47'; Value = '0.07109' },
    @{ Cell = 'D48'; Value = '1.888' },
    @{ Cell = 'D50'; Value = '112.61' },
    @{ Cell = 'D51'; Value = '1.096' }
)

foreach ($u in $priceUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
